# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Insert "Bielorrusia" into the country list ahead of "Ucrania" ---
# In the sheet, data for countries is shown in fixed rows (A44 = Bielorrusia,
# A45 = Ucrania, A46 = Singapur after the update). The underlying data that
# used to belong to "Ucrania" (row 44) now belongs to the newly-surfaced
# "Bielorrusia" row, and the rows below shift down by one: former Ucrania
# data -> row 45 (Ucrania), former Singapur data -> row 46 (Singapur).
# Row 44 receives brand-new "Bielorrusia" totals.

# Row 44 (Bielorrusia) - new data
$ws.Range("A44").Value = "Bielorrusia"
$ws.Range("B44").Value = 4779
$ws.Range("C44").Value = 575
$ws.Range("D44").Value = 342
$ws.Range("E44").Value = 4395
$ws.Range("F44").Value = 65
$ws.Range("G44").Value = 2
$ws.Range("H44").Value = 42

# Row 45 (Ucrania) - takes old Bielorrusia... previously old row-44 data
$ws.Range("A45").Value = "Ucrania"
$ws.Range("B45").Value = 4662
$ws.Range("C45").Value = 501
$ws.Range("D45").Value = 246
$ws.Range("E45").Value = 4291
$ws.Range("F45").Value = 45
$ws.Range("G45").Value = 9
$ws.Range("H45").Value = 125

# Row 46 (Singapur) - takes old row-45 data
$ws.Range("A46").Value = "Singapur"
$ws.Range("B46").Value = 4427
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 683
$ws.Range("E46").Value = 3734
$ws.Range("F46").Value = 29
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 10

# --- Other independent provincia/country updates ---

# Row 20 (Austria)
$ws.Range("B20").Value = 14527
$ws.Range("C20").Value = 51
$ws.Range("D20").Value = 9704
$ws.Range("E20").Value = 4413
$ws.Range("F20").Value = 227

# Row 31 (Rumania)
$ws.Range("B31").Value = 8067
$ws.Range("C31").Value = 360
$ws.Range("D31").Value = 1508
$ws.Range("E31").Value = 6159
$ws.Range("F31").Value = 258

# Row 91 (Libano)
$ws.Range("B91").Value = 668
$ws.Range("C91").Value = 5
$ws.Range("E91").Value = 561
